# DO October 2023 release
# Update the Disease Ontology source_version from v2023-09-29 to v2023-10-21,
# and move the selection from E4 to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Disease Ontology (row 3) source_version: v2023-09-29 -> v2023-10-21
$ws.Range("E3").Value = "v2023-10-21"

# Move the active selection to E3 (was E4)
$ws.Range("E3").Select()
